$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L17").Value = 21437029.5
$ws.Range("N17").Value = -21437365.5
$ws.Range("H17").Value = 6669424.5
$ws.Range("J17").Value = 7145676.5
$ws.Range("K33").Value = 159.72223
$ws.Range("I33").Value = 159.72223
$ws.Range("H33").Value = 153.42105
$ws.Range("M33").Value = 69.27777
$ws.Range("K64").Value = 2751
$ws.Range("H64").Value = 3100.4167
$ws.Range("J64").Value = 3170.3
$ws.Range("M64").Value = -2503
$ws.Range("L64").Value = 3170.3
$ws.Range("I64").Value = 2751
$ws.Range("N64").Value = -3666.3
$ws.Range("M67").Value = -1893
$ws.Range("I67").Value = 2751
$ws.Range("N67").Value = -4886.3
$ws.Range("L67").Value = 3170.3
$ws.Range("K67").Value = 2751
$ws.Range("H67").Value = 3100.4167
$ws.Range("J67").Value = 3170.3
$ws.Range("I98").Value = 767.8095
$ws.Range("J98").Value = 1750
$ws.Range("M98").Value = 730.1905
$ws.Range("N98").Value = -4746
$ws.Range("K98").Value = 767.8095
$ws.Range("H98").Value = 924.96
$ws.Range("L98").Value = 1750
$ws.Range("K100").Value = 1552.1428
$ws.Range("H100").Value = 2298.2942
$ws.Range("I100").Value = 1552.1428
$ws.Range("M100").Value = -1011.1428
$ws.Range("K122").Value = 2303.4285
$ws.Range("H122").Value = 924.96
$ws.Range("M122").Value = 146.5715
$ws.Range("J122").Value = 1750
$ws.Range("N122").Value = -10150
$ws.Range("L122").Value = 5250
$ws.Range("I122").Value = 767.8095
$ws.Range("L124").Value = 37770
$ws.Range("N124").Value = -47590
$ws.Range("H124").Value = 37770
$ws.Range("J124").Value = 37770
$ws.Range("L129").Value = 494745.6900000001
$ws.Range("H129").Value = 139756.9
$ws.Range("N129").Value = -504745.6900000001
$ws.Range("J129").Value = 164915.23
$ws.Range("H132").Value = 2903.4707
$ws.Range("L132").Value = 4500
$ws.Range("N132").Value = -9560
$ws.Range("M132").Value = -6443.5625
$ws.Range("K132").Value = 8973.5625
$ws.Range("I132").Value = 2991.1875
$ws.Range("J132").Value = 1500
$ws.Range("H133").Value = 48828
$ws.Range("J133").Value = 48828
$ws.Range("N133").Value = -58948
$ws.Range("L133").Value = 48828
$ws.Range("J136").Value = 0
$ws.Range("H136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L2").Value = 900
$ws.Range("N2").Value = -1126
$ws.Range("H2").Value = 918.1
$ws.Range("J2").Value = 900
$ws.Range("I32").Value = 9702.143
$ws.Range("M32").Value = -9415.143
$ws.Range("H32").Value = 11414.363
$ws.Range("K32").Value = 9702.143
$ws.Range("I45").Value = 5420.2
$ws.Range("H45").Value = 4992.1113
$ws.Range("M45").Value = -5043.2
$ws.Range("J45").Value = 4457
$ws.Range("K45").Value = 5420.2
$ws.Range("L45").Value = 4457
$ws.Range("N45").Value = -5211
$ws.Range("M61").Value = -1456.3
$ws.Range("L61").Value = 4000
$ws.Range("K61").Value = 1668.3
$ws.Range("H61").Value = 2134.64
$ws.Range("J61").Value = 4000
$ws.Range("I61").Value = 1668.3
$ws.Range("N61").Value = -4424
$ws.Range("J74").Value = 4234.8
$ws.Range("I74").Value = 71429256
$ws.Range("L74").Value = 4234.8
$ws.Range("H74").Value = 41668830
$ws.Range("K74").Value = 71429256
$ws.Range("M74").Value = -71428382
$ws.Range("N74").Value = -5982.8
$ws.Range("K77").Value = 357146280
$ws.Range("N77").Value = -29910
$ws.Range("L77").Value = 21174
$ws.Range("M77").Value = -357141912
$ws.Range("I77").Value = 71429256
$ws.Range("H77").Value = 41668830
$ws.Range("J77").Value = 4234.8
$ws.Range("J97").Value = 333335260
$ws.Range("K97").Value = 1648.7
$ws.Range("M97").Value = -1152.7
$ws.Range("I97").Value = 1648.7
$ws.Range("N97").Value = -333336252
$ws.Range("L97").Value = 333335260
$ws.Range("H97").Value = 76924790
$ws.Range("I102").Value = 1414.5
$ws.Range("K102").Value = 1414.5
$ws.Range("M102").Value = 207.5
$ws.Range("H102").Value = 1531.6
$ws.Range("J110").Value = 1478.75
$ws.Range("L110").Value = 1478.75
$ws.Range("I110").Value = 932.2727
$ws.Range("K110").Value = 932.2727
$ws.Range("M110").Value = 1112.7273
$ws.Range("H110").Value = 1078
$ws.Range("N110").Value = -5568.75
$ws.Range("L116").Value = 900
$ws.Range("J116").Value = 900
$ws.Range("H116").Value = 918.1
$ws.Range("N116").Value = -5488
$ws.Range("M136").Value = -2454.9
$ws.Range("J136").Value = 4000
$ws.Range("H136").Value = 2134.64
$ws.Range("L136").Value = 12000
$ws.Range("I136").Value = 1668.3
$ws.Range("N136").Value = -17100
$ws.Range("K136").Value = 5004.9
$ws.Range("H138").Value = 50237.332
$ws.Range("N138").Value = -60517.332
$ws.Range("L138").Value = 50237.332
$ws.Range("J138").Value = 50237.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 918.1
$ws.Range("N3").Value = -1128
$ws.Range("J3").Value = 900
$ws.Range("L3").Value = 900
$ws.Range("M20").Value = -4878.625
$ws.Range("I20").Value = 5125.625
$ws.Range("H20").Value = 5125.625
$ws.Range("K20").Value = 5125.625
$ws.Range("H105").Value = 1924869.2
$ws.Range("K105").Value = 1599.8572
$ws.Range("M105").Value = 147.1428000000001
$ws.Range("I105").Value = 1599.8572
$ws.Range("H134").Value = 4090.0938
$ws.Range("I134").Value = 4301.1
$ws.Range("J134").Value = 925
$ws.Range("N134").Value = -7845
$ws.Range("M134").Value = -10368.3
$ws.Range("L134").Value = 2775
$ws.Range("K134").Value = 12903.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M39").ClearContents()
$ws.Range("K39").Value = 0
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("H134").Value = 886.65216
$ws.Range("I134").Value = 803.8333
$ws.Range("J134").Value = 1184.8
$ws.Range("N134").Value = -8624.4
$ws.Range("M134").Value = 123.5001000000002
$ws.Range("L134").Value = 3554.4
$ws.Range("K134").Value = 2411.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K122").Value = 3841.71435
$ws.Range("H122").Value = 640.3333
$ws.Range("M122").Value = -1391.71435
$ws.Range("J122").Value = 728.2353000000001
$ws.Range("N122").Value = -11454.1177
$ws.Range("L122").Value = 6554.117700000001
$ws.Range("I122").Value = 426.85715
$ws.Range("H131").Value = 677.38
$ws.Range("J131").Value = 699
$ws.Range("K131").Value = 1286.25
$ws.Range("L131").Value = 2097
$ws.Range("I131").Value = 428.75
$ws.Range("N131").Value = -12177
$ws.Range("M131").Value = 3753.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("H55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("I102").Value = 1774.7826
$ws.Range("K102").Value = 1774.7826
$ws.Range("L102").Value = 3662.8
$ws.Range("J102").Value = 3662.8
$ws.Range("N102").Value = -6906.8
$ws.Range("M102").Value = -152.7826
$ws.Range("H102").Value = 2111.9285
$ws.Range("L106").Value = 10000
$ws.Range("N106").Value = -12524
$ws.Range("J106").Value = 10000
$ws.Range("H106").Value = 10000
$ws.Range("K122").Value = 6000
$ws.Range("H122").Value = 7000
$ws.Range("M122").Value = -3550
$ws.Range("I122").Value = 2000
$ws.Range("H132").Value = 20910.967
$ws.Range("L132").Value = 161044.5
$ws.Range("N132").Value = -166104.5
$ws.Range("M132").Value = -11047.1
$ws.Range("K132").Value = 13577.1
$ws.Range("I132").Value = 4525.7
$ws.Range("J132").Value = 53681.5
$ws.Range("H138").Value = 45000
$ws.Range("N138").Value = -55280
$ws.Range("L138").Value = 45000
$ws.Range("J138").Value = 45000
$ws.Range("J141").Value = 54811
$ws.Range("L141").Value = 54811
$ws.Range("H141").Value = 54811
$ws.Range("N141").Value = -65171

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 2588
$ws.Range("J24").Value = 2588
$ws.Range("L24").Value = 2588
$ws.Range("N24").Value = -3274
$ws.Range("I46").Value = 4158.3335
$ws.Range("M46").Value = -3970.3335
$ws.Range("K46").Value = 4158.3335
$ws.Range("H46").Value = 3465.2
$ws.Range("J55").Value = 217
$ws.Range("L55").Value = 217
$ws.Range("I55").Value = 1372.8572
$ws.Range("M55").Value = -1199.8572
$ws.Range("H55").Value = 891.25
$ws.Range("N55").Value = -563
$ws.Range("K55").Value = 1372.8572
$ws.Range("H130").Value = 16625
$ws.Range("J130").Value = 16625
$ws.Range("N130").Value = -26665
$ws.Range("L130").Value = 16625
$ws.Range("M136").Value = -1786.7307
$ws.Range("J136").Value = 3159
$ws.Range("H136").Value = 1766.8438
$ws.Range("L136").Value = 9477
$ws.Range("I136").Value = 1445.5769
$ws.Range("N136").Value = -14577
$ws.Range("K136").Value = 4336.7307

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J74").Value = 50863
$ws.Range("L74").Value = 50863
$ws.Range("H74").Value = 50863
$ws.Range("N74").Value = -52735
$ws.Range("N77").Value = -161949
$ws.Range("L77").Value = 152589
$ws.Range("H77").Value = 50863
$ws.Range("J77").Value = 50863
$ws.Range("H105").Value = 29399.666
$ws.Range("J105").Value = 29399.666
$ws.Range("L105").Value = 29399.666
$ws.Range("N105").Value = -36387.666
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("M136").Value = -88480038
$ws.Range("J136").Value = 8866.666999999999
$ws.Range("H136").Value = 23463106
$ws.Range("L136").Value = 26600.001
$ws.Range("I136").Value = 29494196
$ws.Range("N136").Value = -31700.001
$ws.Range("K136").Value = 88482588
